$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Idioma" column header
$ws.Range("E1").Value2 = "Idioma"

# New language values for existing rows
$ws.Range("E2").Value2 = "Japones"
$ws.Range("E3").Value2 = "Japones"

# New 4th data row (repeats row2's A/B/C/D, new language Mandarim)
$ws.Range("A4").Value2 = $ws.Range("A2").Value2
$ws.Range("B4").Value2 = $ws.Range("B2").Value2
$ws.Range("C4").Value2 = $ws.Range("C2").Value2
$ws.Range("D4").Value2 = $ws.Range("D2").Value2
$ws.Range("E4").Value2 = "Mandarim"

# Autofit the new column to its content
$ws.Range("E1:E4").EntireColumn.AutoFit()

# New empty, underlined cell further down the sheet
$ws.Range("F8").Font.Underline = $true

# Move / leave selection on the last-touched cell
$ws.Range("F8").Select() | Out-Null
